$wb = $excel.ActiveWorkbook

# Locate the last existing sheet (2025-09-08) to insert the new sheet after it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-09-15"

# Header row
$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "author"
$newSheet.Range("D1").Value = "latest_episode"

$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Ranking data rows (rank, title, author, latest_episode)
$data = @(
    @(1, "宇崎ちゃんは遊びたい！", "丈(著者)", "第127話"),
    @(2, "生徒会にも穴はある！", "むちまろ", "第135話`tムラ塚先生の憂鬱"),
    @(3, "不徳のギルド", "河添太一", "第９８話：コインの裏表"),
    @(4, "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第６９話『岩鬼停止』②"),
    @(5, "魔術師クノンは見えている", "La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)", "第40話②"),
    @(6, "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第32話：思考を奪う②"),
    @(7, "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第64話 更衣室の戦いが始まった（２）"),
    @(8, "実は俺、最強でした？", "原作：澄守 彩 漫画：高橋 愛", "第125話　四騎戦の打ち上げです！"),
    @(9, "いとこのこ", "いぬちく(著者)", "第38話"),
    @(10, "追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。", "六志麻あさ 業務用餅 kisui", "第７１話"),
    @(11, "クラスで２番目に可愛い女の子と友だちになった", "尾野凛(漫画) たかた(原作) 日向あずり(キャラクター原案)", "第35話①"),
    @(12, "男女比1：5の世界でも普通に生きられると思った？　～激重感情な彼女たちが無自覚男子に翻弄されたら～", "三藤 孝太郎(原作) 桃季憂(漫画) jimmy(キャラクター原案)", "第11話-1"),
    @(13, "美人女上司滝沢さん", "やんBARU(著者)", "第203話"),
    @(14, "蜘蛛ですが、なにか？", "かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)", "第76話その1"),
    @(15, "勇者に全部奪われた俺は勇者の母親とパーティを組みました！", "久遠まこと(著者) 石のやっさん(原作)", "第30話"),
    @(16, "ダンジョンの幼なじみ", "久真やすひさ(著者)", "第57話"),
    @(17, "辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？", "tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)", "番外編①"),
    @(18, "異世界魔王と召喚少女の奴隷魔術", "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大", "第128話　レムと話してみる（中編）"),
    @(19, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５２話　暴走を止める器用貧乏（３）"),
    @(20, "不純な彼女達は懺悔しない", "ポロロッカ(著者)", "第31話"),
    @(21, "異世界のんびり農家", "剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)", "第304話"),
    @(22, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "第２３食　巨大ヘビモンスターさん、パクパクですわ！（４）"),
    @(23, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第139話 よくわからないけれどズゥゥゥンときたみたいです（１）"),
    @(24, "望まぬ不死の冒険者", "中曽根ハイジ（漫画） 丘野 優（原作） じゃいあん（キャラクター原案）", "第60話　異国の観光客"),
    @(25, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第73話①"),
    @(26, "魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～", "漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ", "第47話 魔導具師とつながれたもの④"),
    @(27, "「ククク……。奴は四天王の中でも最弱」と解雇された俺、なぜか勇者と聖女の師匠になる", "漫画：芳橋アツシ 原作：延野正行 キャラクター原案：坂野杏梨", "第41話　奴は帰りたい（後編）"),
    @(28, "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜", "戸賀 環 坂木持丸 riritto", "第52話②　最強の武器を手に入れてみた"),
    @(29, "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～", "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき", "第33話 独身貴族は見積もりを誤る（4）"),
    @(30, "解雇された暗黒兵士(30代)のスローなセカンドライフ", "岡沢六十四 るれくちぇ sage・ジョー", "第72話(後編) ダリエルの照らした道"),
    @(31, "アイドル辞めるけど結婚してくれますか!?", "三吉汐美(著者)", "第17話前半"),
    @(32, "くらいあの子としたいこと", "碇マナツ(著者)", "第84話"),
    @(33, "姫様“拷問”の時間です", "原作:春原ロビンソン　漫画:ひらけい", "拷問149"),
    @(34, "治癒魔法の間違った使い方 ~戦場を駆ける回復要員~", "九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)", "第82話その3"),
    @(35, "落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～", "村上よしゆき 茨木野 あるてら", "第４２話　勇者、六邪神将相手に舐めプしてたら、ピンチになる（２）"),
    @(36, "サーシャちゃんとクラスメイトオタクくん", "はぐはぐ(著者)", "第85話"),
    @(37, "十年目、帰還を諦めた転移者はいまさら主人公になる", "原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう", "第１９話①"),
    @(38, "理想のヒモ生活", "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)", "第87話　その2"),
    @(39, "バーサス", "原作：ONE 漫画：あずま京太郎 構成：bose", "第28話　工場（1）"),
    @(40, "三枝さんはメガネ先輩と恋を描く", "セレビィ量産型(著者)", "第22話前半"),
    @(41, "塔の管理をしてみよう", "盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)", "第92話後編"),
    @(42, "ライブダンジョン！", "ことりりょう(作画) dy冷凍(原作) Mika Pikazo(キャラクター原案)", "第89話前半"),
    @(43, "配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信", "下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)", "第2話後編"),
    @(44, "異世界でも無難に生きたい症候群", "原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう", "第31話③"),
    @(45, "魔王になったので、ダンジョン造って人外娘とほのぼのする", "遠野ノオト(作画) 流優(原作) だぶ竜(キャラクター原案)", "第1話後半"),
    @(46, "女子高生の無駄づかい", "ビーノ(著者)", "第135話　きずな"),
    @(47, "路地裏で拾った女の子がバッドエンド後の乙女ゲームのヒロインだった件", "カボチャマスク(原作) 樋乃えなが(作画) へいろー(キャラクター原案)", "第2話"),
    @(48, "ゲーム世界で魔物に転生してしまった俺、前世で推しだったヒロインを拾ってしまう", "三部べべ(漫画) ねうしとら(原作)", "第1話"),
    @(49, "隣の席のヤンキー清水さんが髪を黒く染めてきた", "底花(原作) 真田若楓(漫画) ハム(キャラクター原案)", "第11話-1"),
    @(50, "最弱貴族に転生したので悪役たちを集めてみた", "空野進 sorani ファルまろ", "第12話　最弱貴族、悪役令嬢と出会う（２）")
)

$rowCount = $data.Count
$colCount = 4
$arr = New-Object 'object[,]' $rowCount, $colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$targetRange = $newSheet.Range("A2").Resize($rowCount, $colCount)
$targetRange.Value = $arr

Write-Host "Added sheet 2025-09-15 with $rowCount data rows"
